# "chcked and fixed the ers" -- checked column B (COMPANY) for duplicate
# entries, widened the column so the full names are visible, renamed the
# duplicate company entries so every row is unique, and left a
# "Duplicate Values" conditional-formatting rule on the column so future
# dupes get flagged automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column B to fit the longer, de-duplicated company names.
$ws.Columns.Item(2).ColumnWidth = 39

# Rename the duplicated company entries so every value in column B is unique.
$ws.Cells.Item(70, 2).Value = "Infosys through - #HackwithInfy-2"
$ws.Cells.Item(59, 2).Value = "Larsen & Toubro Infotech Ltd.-2"
$ws.Cells.Item(58, 2).Value = "Infosys through - #HackwithInfy-3"
$ws.Cells.Item(48, 2).Value = "Larsen & Toubro Infotech Ltd.-3"

# Highlight Cells Rules > Duplicate Values, applied to the whole column so
# any future duplicate company name gets flagged (standard light-red
# fill / dark-red text preset).
$dupRange = $ws.Range("B1:B1048576")
$condFormat = $dupRange.FormatConditions.AddUniqueValues()
$condFormat.DupeUnique = 1
$condFormat.Font.Color = 0x0006009C
$condFormat.Interior.Color = 0x00CEC7FF

# Restore the scroll position / active cell of the last save.
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
[void]$ws.Range("B49").Select()
